$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row at position 38 (new claim: -416 / Paraguay 3765)
$ws.Rows(38).Insert()

$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = '-416'
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = '5/26/2025'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'Paraguay 3765'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '806926557'
$ws.Range("F38").NumberFormat = "@"
$ws.Range("F38").Value = 'Optical Power'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = 'Pendiente'
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = 'Colocar columna para pedir traspaso de nodo teco'
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = '1'
$ws.Range("J38").NumberFormat = "@"
$ws.Range("J38").Value = 'Cambio'
$ws.Range("K38").NumberFormat = "@"
$ws.Range("K38").Value = 'Nodo Teco'
$ws.Range("L38").NumberFormat = "@"
$ws.Range("L38").Value = 'Pasante'
$ws.Range("M38").Value = -58.416562
$ws.Range("N38").Value = -34.590589

# Insert new row at position 65 (new claim: 6076 / MATHEU 727)
$ws.Rows(65).Insert()

$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = '6076'
$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = '6/24/2025'
$ws.Range("C65").NumberFormat = "@"
$ws.Range("C65").Value = 'MATHEU 727'
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = '3'
$ws.Range("E65").NumberFormat = "@"
$ws.Range("E65").Value = '807763063'
$ws.Range("F65").NumberFormat = "@"
$ws.Range("F65").Value = 'Optical Power'
$ws.Range("G65").NumberFormat = "@"
$ws.Range("G65").Value = 'Pendiente'
$ws.Range("H65").NumberFormat = "@"
$ws.Range("H65").Value = 'Colocar R400 para pedir a base traspaso de nodo propio y posterior a TLC'
$ws.Range("I65").NumberFormat = "@"
$ws.Range("I65").Value = '1'
$ws.Range("J65").NumberFormat = "@"
$ws.Range("J65").Value = 'Cambio'
$ws.Range("K65").NumberFormat = "@"
$ws.Range("K65").Value = 'Nodo TLC'
$ws.Range("L65").NumberFormat = "@"
$ws.Range("L65").Value = 'Pasante'
$ws.Range("M65").Value = -58.400169
$ws.Range("N65").Value = -34.617784
